# Week 5 update:
# For every match row (rows 2-66): the "(N)" seed/score suffix on both
# team names (columns A & B) is reset to "(0)", the Drawn flag (column C)
# is set to 13, and the Winning Team (column D) becomes "Not Determined"
# since the result could not be determined yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: when assigning a string to a cell's .Value, Excel treats a
# leading single-quote as a "force text" entry marker and swallows it
# rather than storing it (and marks the cell with a quote-prefix style).
# Double it up so a literal leading quote survives the round trip, then
# restore the cell's original (unflagged) style.
function Set-CellText($cell, [string]$text) {
    if ($text.StartsWith("'")) {
        $cell.Value = "'" + $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

$lastRow = 66
for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)
    $cellC = $ws.Cells.Item($r, 3)
    $cellD = $ws.Cells.Item($r, 4)

    $teamA = $cellA.Value2
    $teamB = $cellB.Value2

    $newTeamA = [System.Text.RegularExpressions.Regex]::Replace($teamA, '\(\d+\)$', '(0)')
    $newTeamB = [System.Text.RegularExpressions.Regex]::Replace($teamB, '\(\d+\)$', '(0)')

    Set-CellText $cellA $newTeamA
    Set-CellText $cellB $newTeamB
    $cellC.Value = 13
    $cellD.Value = "Not Determined"
}
